$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24:B24").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)

$ws.Range("A25").Value = "LOME"
$ws.Range("B25").Value = "x"

$ws.Range("D23").Select()
